# "Repayment schedule" is the active sheet in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new (blank) column was inserted at column N, pushing the former
# N/O/P columns (Late / Date(heading) / Outstanding + their data) one
# slot to the right, to O/P/Q respectively.
$ws.Columns("N:N").Insert()

# Give the newly inserted column a width close to its left neighbour
# (column M), matching how the new column N ends up sized in the
# target workbook.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Reflect the new selection/active cell left after the edit.
$ws.Range("Q10").Select()
